# Fix mojibake: "Â±" (UTF-8 bytes for ± mis-decoded as Latin-1) -> "±"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = "$([string]([char]0x00C2))$([string]([char]0x00B1))"
$goodChar = [string]([char]0x00B1)

$range = $ws.Range("B2:D17")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains($badChar)) {
        $cell.Value2 = $val.Replace($badChar, $goodChar)
    }
}
